$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows (columns B and D) with refreshed source data
$ws.Cells.Item(205, 2).Value = 4373594970000
$ws.Cells.Item(205, 4).Value = 195008452486.4765
$ws.Cells.Item(206, 2).Value = 4411934620000
$ws.Cells.Item(206, 4).Value = 196098882059.7977
$ws.Cells.Item(210, 2).Value = 4577407590000
$ws.Cells.Item(210, 4).Value = 197829026890.5965
$ws.Cells.Item(211, 2).Value = 4566459490000
$ws.Cells.Item(211, 4).Value = 204213160361.4111
$ws.Cells.Item(212, 2).Value = 4592275590000
$ws.Cells.Item(212, 4).Value = 198028441162.9506
$ws.Cells.Item(213, 2).Value = 4639859400000
$ws.Cells.Item(213, 4).Value = 196360924007.7125
$ws.Cells.Item(214, 2).Value = 4641345140000
$ws.Cells.Item(214, 4).Value = 196115881881.1865
$ws.Cells.Item(215, 2).Value = 4681223420000
$ws.Cells.Item(215, 4).Value = 204604827761.64
$ws.Cells.Item(216, 2).Value = 4725508480000
$ws.Cells.Item(216, 4).Value = 203831027810.3441
$ws.Cells.Item(217, 2).Value = 4680322510000
$ws.Cells.Item(217, 4).Value = 206432055012.9981
$ws.Cells.Item(219, 2).Value = 4809150480000
$ws.Cells.Item(219, 4).Value = 208310072120.0702
$ws.Cells.Item(221, 2).Value = 4958595660000
$ws.Cells.Item(221, 4).Value = 200173330572.0967
$ws.Cells.Item(222, 2).Value = 5004666910000
$ws.Cells.Item(222, 4).Value = 206283339234.405
$ws.Cells.Item(223, 2).Value = 5020790900000
$ws.Cells.Item(223, 4).Value = 211475627271.8466
$ws.Cells.Item(224, 2).Value = 5059232680000
$ws.Cells.Item(224, 4).Value = 227082617645.2899
$ws.Cells.Item(225, 2).Value = 5094308060000
$ws.Cells.Item(225, 4).Value = 231533988533.9885
$ws.Cells.Item(226, 2).Value = 5178041490000
$ws.Cells.Item(226, 4).Value = 224145123064.5561
$ws.Cells.Item(227, 2).Value = 5214187690000
$ws.Cells.Item(227, 4).Value = 222760943901.6367
$ws.Cells.Item(228, 2).Value = 5235568230000
$ws.Cells.Item(228, 4).Value = 237994978328.2156
$ws.Cells.Item(229, 2).Value = 5179738620000
$ws.Cells.Item(229, 4).Value = 241084591915.8116
$ws.Cells.Item(230, 2).Value = 5290478980000
$ws.Cells.Item(230, 4).Value = 246567659217.4866
$ws.Cells.Item(231, 2).Value = 5390398340000
$ws.Cells.Item(231, 4).Value = 248749920512.0455
$ws.Cells.Item(232, 2).Value = 5449356120000
$ws.Cells.Item(232, 4).Value = 244698227636.6543
$ws.Cells.Item(233, 2).Value = 5471474170000
$ws.Cells.Item(233, 4).Value = 254209314005.6218
$ws.Cells.Item(234, 2).Value = 5507491430000
$ws.Cells.Item(234, 4).Value = 264768603178.1897
$ws.Cells.Item(236, 2).Value = 5564521500000
$ws.Cells.Item(236, 4).Value = 259066790507.9822
$ws.Cells.Item(239, 2).Value = 5617130550000
$ws.Cells.Item(239, 4).Value = 253143058331.5119
$ws.Cells.Item(240, 2).Value = 5647837280000
$ws.Cells.Item(240, 4).Value = 250967585233.9003
$ws.Cells.Item(241, 2).Value = 5542014840000
$ws.Cells.Item(241, 4).Value = 253256051601.582
$ws.Cells.Item(242, 2).Value = 5630383690000
$ws.Cells.Item(242, 4).Value = 259884453460.5434
$ws.Cells.Item(243, 2).Value = 5704249840000
$ws.Cells.Item(243, 4).Value = 254295839565.2552
$ws.Cells.Item(244, 2).Value = 5739159050000
$ws.Cells.Item(244, 4).Value = 260129496843.0881
$ws.Cells.Item(245, 2).Value = 5742427260000
$ws.Cells.Item(245, 4).Value = 246010541422.8308
$ws.Cells.Item(246, 2).Value = 5825723830000
$ws.Cells.Item(246, 4).Value = 252980514452.8375
$ws.Cells.Item(247, 2).Value = 5801917230000
$ws.Cells.Item(247, 4).Value = 245858014335.5658
$ws.Cells.Item(248, 2).Value = 5833040250000
$ws.Cells.Item(248, 4).Value = 242345086968.173
$ws.Cells.Item(249, 2).Value = 5855415460000
$ws.Cells.Item(249, 4).Value = 240184235925.8127
$ws.Cells.Item(250, 2).Value = 5887405600000
$ws.Cells.Item(250, 4).Value = 234668314187.5463
$ws.Cells.Item(251, 2).Value = 5915934540000
$ws.Cells.Item(251, 4).Value = 238841410132.7046
$ws.Cells.Item(252, 2).Value = 5940210650000
$ws.Cells.Item(252, 4).Value = 253955001539.0667
$ws.Cells.Item(256, 2).Value = 6077524080000
$ws.Cells.Item(256, 4).Value = 280799406754.406
$ws.Cells.Item(258, 2).Value = 6224248910000
$ws.Cells.Item(258, 4).Value = 280487810714.3501

# Append new month row 259, copying formatting (incl. date style) from row 258
$ws.Cells.Item(258, 1).Copy($ws.Cells.Item(259, 1))
$ws.Cells.Item(259, 1).Value = 45078
$ws.Cells.Item(259, 2).Value = 6243183470000
$ws.Cells.Item(259, 3).Value = 0.04592443593311565
$ws.Cells.Item(259, 4).Value = 286714679286.7017
